# EWD-23048 - External review: Read comments in the app
#
# - "reviewTabTitle" row: EN text changes from "Share course for review" to "Comments"
# - NL/DE translations (columns D/E) for the "Review" section (rows 351-357) are cleared
#   (no longer translated / placeholders removed)
# - A new localization row is added for the key "reviewNoComments" = "No comments yet"
#   right after the "reviewPublishingCourse" row, pushing the "Feedback" block down by one row
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the English text for reviewTabTitle (row 352, column C = EN)
$ws.Range("C352").Value = "Comments"

# Remove the NL (D) and EN-mirrored (E) translations for the Review section rows 351-357
$ws.Range("D351:E357").ClearContents()

# Insert a new row for the "reviewNoComments" key right after row 357 (reviewPublishingCourse),
# shifting every row from 358 onward down by one
$ws.Rows.Item(358).Insert()
$ws.Range("B358").Value = "reviewNoComments"
$ws.Range("C358").Value = "No comments yet"
